{"js": "// The document contains a single 5x3 table of lattice-multiplication\n// exercises. Each cell holds one paragraph / one run with several lines\n// (separated by manual line breaks <w:br/>, i.e. the vertical-tab char\n// \"\\u000b\" in the Word.js text model):\n//   \"AA x BB\"\n//   \"  b1    b2\"\n//   \"  ----\"\n//   \"a1|    |\"\n//   \"a2|    |\"\n// This edit swaps each cell's two factors (AA x BB) for a new pair,\n// keeping every other aspect (run formatting, paragraph/table layout)\n// unchanged. Build the 15 replacement strings up front, then push each\n// one into its cell in row-major order.\n\nconst newCells = [\n  \"76 x 76\\u000b  7    6\\u000b  ----\\u000b7|    |\\u000b6|    |\",\n  \"62 x 48\\u000b  4    8\\u000b  ----\\u000b6|    |\\u000b2|    |\",\n  \"41 x 49\\u000b  4    9\\u000b  ----\\u000b4|    |\\u000b1|    |\",\n  \"53 x 49\\u000b  4    9\\u000b  ----\\u000b5|    |\\u000b3|    |\",\n  \"13 x 83\\u000b  8    3\\u000b  ----\\u000b1|    |\\u000b3|    |\",\n  \"21 x 67\\u000b  6    7\\u000b  ----\\u000b2|    |\\u000b1|    |\",\n  \"94 x 22\\u000b  2    2\\u000b  ----\\u000b9|    |\\u000b4|    |\",\n  \"35 x 37\\u000b  3    7\\u000b  ----\\u000b3|    |\\u000b5|    |\",\n  \"53 x 21\\u000b  2    1\\u000b  ----\\u000b5|    |\\u000b3|    |\",\n  \"29 x 47\\u000b  4    7\\u000b  ----\\u000b2|    |\\u000b9|    |\",\n  \"59 x 31\\u000b  3    1\\u000b  ----\\u000b5|    |\\u000b9|    |\",\n  \"16 x 70\\u000b  7    0\\u000b  ----\\u000b1|    |\\u000b6|    |\",\n  \"65 x 36\\u000b  3    6\\u000b  ----\\u000b6|    |\\u000b5|    |\",\n  \"33 x 69\\u000b  6    9\\u000b  ----\\u000b3|    |\\u000b3|    |\",\n  \"86 x 42\\u000b  4    2\\u000b  ----\\u000b8|    |\\u000b6|    |\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 3;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n    range.insertText(newCells[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 5x3 table of lattice-multiplication\n# exercises. Each cell holds one paragraph / one run with several lines\n# (separated by manual line breaks, vertical-tab char 0x0B):\n#   \"AA x BB\"\n#   \"  b1    b2\"\n#   \"  ----\"\n#   \"a1|    |\"\n#   \"a2|    |\"\n# This edit swaps each cell's two factors (AA x BB) for a new pair,\n# keeping every other aspect (run formatting, paragraph/table layout)\n# unchanged. Build the 15 replacement strings (row-major order) and push\n# each one into its cell, excluding the trailing end-of-cell mark.\n\n$d = $word.ActiveDocument\n\n$vt = [char]11\n\n$newCells = @(\n    (\"76 x 76\" + $vt + \"  7    6\" + $vt + \"  ----\" + $vt + \"7|    |\" + $vt + \"6|    |\"),\n    (\"62 x 48\" + $vt + \"  4    8\" + $vt + \"  ----\" + $vt + \"6|    |\" + $vt + \"2|    |\"),\n    (\"41 x 49\" + $vt + \"  4    9\" + $vt + \"  ----\" + $vt + \"4|    |\" + $vt + \"1|    |\"),\n    (\"53 x 49\" + $vt + \"  4    9\" + $vt + \"  ----\" + $vt + \"5|    |\" + $vt + \"3|    |\"),\n    (\"13 x 83\" + $vt + \"  8    3\" + $vt + \"  ----\" + $vt + \"1|    |\" + $vt + \"3|    |\"),\n    (\"21 x 67\" + $vt + \"  6    7\" + $vt + \"  ----\" + $vt + \"2|    |\" + $vt + \"1|    |\"),\n    (\"94 x 22\" + $vt + \"  2    2\" + $vt + \"  ----\" + $vt + \"9|    |\" + $vt + \"4|    |\"),\n    (\"35 x 37\" + $vt + \"  3    7\" + $vt + \"  ----\" + $vt + \"3|    |\" + $vt + \"5|    |\"),\n    (\"53 x 21\" + $vt + \"  2    1\" + $vt + \"  ----\" + $vt + \"5|    |\" + $vt + \"3|    |\"),\n    (\"29 x 47\" + $vt + \"  4    7\" + $vt + \"  ----\" + $vt + \"2|    |\" + $vt + \"9|    |\"),\n    (\"59 x 31\" + $vt + \"  3    1\" + $vt + \"  ----\" + $vt + \"5|    |\" + $vt + \"9|    |\"),\n    (\"16 x 70\" + $vt + \"  7    0\" + $vt + \"  ----\" + $vt + \"1|    |\" + $vt + \"6|    |\"),\n    (\"65 x 36\" + $vt + \"  3    6\" + $vt + \"  ----\" + $vt + \"6|    |\" + $vt + \"5|    |\"),\n    (\"33 x 69\" + $vt + \"  6    9\" + $vt + \"  ----\" + $vt + \"3|    |\" + $vt + \"3|    |\"),\n    (\"86 x 42\" + $vt + \"  4    2\" + $vt + \"  ----\" + $vt + \"8|    |\" + $vt + \"6|    |\")\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $rng = $cell.Range\n        # Drop the trailing end-of-cell marker (CR + cell-mark) so we\n        # only overwrite the visible content, same as Word does when you\n        # select-and-type inside a cell.\n        $rng.End = $rng.End - 1\n        $rng.Text = $newCells[$idx]\n        $idx++\n    }\n}\n"}
